$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.442.86"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "2.609.35"
$ws.Range("E3").Value = "  +9.59%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.87"
$ws.Range("E5").Value = "  +1.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.25"
$ws.Range("E6").Value = "  +1.98%  "
$ws.Range("E7").Value = "  +5.76%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +11.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.95"
$ws.Range("E10").Value = "  +12.32%  "
$ws.Range("E11").Value = "  +5.66%  "
$ws.Range("E12").Value = "  +14.58%  "
$ws.Range("D13").Value = "3.007.46"
$ws.Range("E13").Value = "  +9.68%  "
$ws.Range("D15").Value = "2.630.74"
$ws.Range("E15").Value = "  +10.83%  "
$ws.Range("E16").Value = "  +10.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.89"
$ws.Range("E17").Value = "  +8.68%  "
$ws.Range("D18").Value = "46.595.66"
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.33"
$ws.Range("E19").Value = "  +2.95%  "
$ws.Range("E20").Value = "  +5.03%  "
$ws.Range("E21").Value = "  +10.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.08"
$ws.Range("E22").Value = "  +5.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "257.30"
$ws.Range("E23").Value = "  +4.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.02"
$ws.Range("E24").Value = "  +7.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.25"
$ws.Range("E25").Value = "  +15.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.46"
$ws.Range("E26").Value = "  +34.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  +6.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.30"
$ws.Range("E29").Value = "  +0.65%  "
$ws.Range("E30").Value = "  +3.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.20"
$ws.Range("E31").Value = "  +11.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.73"
$ws.Range("E32").Value = "  -2.08%  "
$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.98"
$ws.Range("E33").Value = "  +6.00%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.33"
$ws.Range("E34").Value = "  +21.88%  "
$ws.Range("E35").Value = "  +7.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "150.81"
$ws.Range("E36").Value = "  +2.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.118"
$ws.Range("E37").Value = "  +4.22%  "
$ws.Range("E38").Value = "  +5.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.86"
$ws.Range("E39").Value = "  +5.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.18"
$ws.Range("E40").Value = "  +5.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.63"
$ws.Range("E41").Value = "  +12.37%  "
$ws.Range("E42").Value = "  +7.46%  "
$ws.Range("D43").Value = "2.046.80"
$ws.Range("E43").Value = "  +6.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.15"
$ws.Range("E44").Value = "  +33.87%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.66"
$ws.Range("E46").Value = "  -0.83%  "
$ws.Range("E47").Value = "  +10.18%  "
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "109.90"
$ws.Range("E49").Value = "  +11.63%  "
$ws.Range("E50").Value = "  +6.98%  "
$ws.Range("D51").Value = "2.864.01"
$ws.Range("E51").Value = "  +9.61%  "
